# Update "想去人数" (interest count) figures in column F on the
# "展览" and "全部类型" worksheets, reflecting the latest scrape
# (gh-pages output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 56
$wsExhibit.Range("F4").Value = 7351
$wsExhibit.Range("F5").Value = 273
$wsExhibit.Range("F6").Value = 436
$wsExhibit.Range("F7").Value = 3877
$wsExhibit.Range("F9").Value = 552
$wsExhibit.Range("F11").Value = 618

# --- Sheet "全部类型" ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 56
$wsAll.Range("F5").Value = 7351
$wsAll.Range("F7").Value = 273
$wsAll.Range("F8").Value = 436
$wsAll.Range("F9").Value = 3877
$wsAll.Range("F11").Value = 552
$wsAll.Range("F13").Value = 618
